$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" -----------
# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

# zh-cn sheet: Status column C, rows 2-3.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

# de-de sheet: Status column C, rows 2-3.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- Shrink the now-narrower "status" columns to fit the new text ----------
# The status text got shorter ("Ready for handoff" -> "In Translation"), so
# the report regenerate also narrowed these columns.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
